# The source workbook tracks weekly Zapallo italiano price observations.
# A new weekly record was inserted before the existing row 349, pushing every
# subsequent row down by one (the former last row, 426, becomes row 427).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 349 - this shifts rows 349:426 down to 350:427 and
# grows the sheet's used range to A1:R427, matching the target dimension.
$ws.Rows.Item(349).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A349").Value = 8
$ws.Range("B349").Value = "Terminal La Palmera de La Serena"
$ws.Range("C349").Value = "Coquimbo"
$ws.Range("D349").Value = 44889
$ws.Range("E349").Value = 4
$ws.Range("F349").Value = 100112032
$ws.Range("G349").Value = "Zapallo italiano"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 460
$ws.Range("K349").Value = 7000
$ws.Range("L349").Value = 8000
$ws.Range("M349").Value = 7500
$ws.Range("N349").Value = "$/caja 60 unidades"
$ws.Range("O349").Value = "Provincia de Limarí"
$ws.Range("P349").Value = 125
$ws.Range("Q349").Value = 60
$ws.Range("R349").Value = "Hortaliza"
